$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.354.36"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "2.299.52"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.96"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.79"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  +2.73%  "
$ws.Range("D9").Value = "2.298.67"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1000"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.56"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "2.704.83"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "58.245.15"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "2.278.45"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("E20").Value = "  -4.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.54"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.60"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.05"
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.27"
$ws.Range("E28").Value = "  -5.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.33"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "0.0₃0722"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.76"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.05"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.79"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.49"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.02"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "289.93"
$ws.Range("E43").Value = "  -4.79%  "
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0951"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0495"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.557"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.15"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.95"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  -0.65%  "
